$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ponds")

# Fix the typo "Amoun (Gal)" -> "Amount (Gal)" in the header cell G3
$ws.Range("G3").Value = "Amount (Gal)"

# Update the selected cell to match the new active cell G3
$ws.Activate()
$ws.Range("G3").Select()
